# 自动更新Excel文件 - daily refresh of remaining-day counters
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $d = $dCell.Value2
    $e = $eCell.Value2
    $f = $fCell.Value2

    if ($null -eq $d -or $null -eq $e -or $null -eq $f) {
        continue
    }

    if ($e -eq $d) {
        # Row has not started its countdown yet (or is malformed) - skip it.
        continue
    }
    elseif ($e -eq 1) {
        # Cycle is about to expire: roll the start date forward to the
        # previous end date and reset the remaining-days counter to the
        # full duration.
        $fStr = [string][int]$f
        $year = [int]$fStr.Substring(0, 4)
        $month = [int]$fStr.Substring(4, 2)
        $day = [int]$fStr.Substring(6, 2)
        $startDate = Get-Date -Year $year -Month $month -Day $day
        $newDate = $startDate.AddDays([double]$d)
        $newF = [int]($newDate.ToString("yyyyMMdd"))

        $fCell.Value = $newF
        $eCell.Value = $d
    }
    else {
        $eCell.Value = $e - 1
    }
}
